# Auto-generated script: update Leve price/profit data cells per scheduled runner refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 309.92307
$ws.Range("I4").Value = 202.63637
$ws.Range("K4").Value = 202.63637
$ws.Range("M4").Value = -88.63637
$ws.Range("H107").Value = 1399.5483
$ws.Range("I107").Value = 1381.5238
$ws.Range("K107").Value = 1381.5238
$ws.Range("M107").Value = 538.4762000000001
$ws.Range("H112").Value = 3592.4856
$ws.Range("J112").Value = 3604.2686
$ws.Range("L112").Value = 10812.8058
$ws.Range("N112").Value = -13028.8058
$ws.Range("H118").Value = 2766.7144
$ws.Range("I118").Value = 2766.7144
$ws.Range("K118").Value = 8300.143199999999
$ws.Range("M118").Value = -6643.143199999999
$ws.Range("H138").Value = 118389.55
$ws.Range("I138").Value = 6054
$ws.Range("K138").Value = 18162
$ws.Range("M138").Value = -13022

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25411.305
$ws.Range("I32").Value = 21284.193
$ws.Range("J32").Value = 50999.4
$ws.Range("K32").Value = 21284.193
$ws.Range("L32").Value = 50999.4
$ws.Range("M32").Value = -20997.193
$ws.Range("N32").Value = -51573.4
$ws.Range("H61").Value = 1863127.5
$ws.Range("I61").Value = 2501730.2
$ws.Range("K61").Value = 2501730.2
$ws.Range("M61").Value = -2501518.2
$ws.Range("H74").Value = 3863295
$ws.Range("I74").Value = 4465847.5
$ws.Range("J74").Value = 6960
$ws.Range("K74").Value = 4465847.5
$ws.Range("L74").Value = 6960
$ws.Range("M74").Value = -4464973.5
$ws.Range("N74").Value = -8708
$ws.Range("H77").Value = 3863295
$ws.Range("I77").Value = 4465847.5
$ws.Range("J77").Value = 6960
$ws.Range("K77").Value = 22329237.5
$ws.Range("L77").Value = 34800
$ws.Range("M77").Value = -22324869.5
$ws.Range("N77").Value = -43536
$ws.Range("H102").Value = 8336921
$ws.Range("I102").Value = 3805.45
$ws.Range("K102").Value = 3805.45
$ws.Range("M102").Value = -2183.45
$ws.Range("H108").Value = 90000
$ws.Range("J108").Value = 90000
$ws.Range("L108").Value = 90000
$ws.Range("N108").Value = -97680
$ws.Range("H110").Value = 1184.7297
$ws.Range("I110").Value = 1030.3549
$ws.Range("K110").Value = 1030.3549
$ws.Range("M110").Value = 1014.6451
$ws.Range("H122").Value = 6305.926
$ws.Range("I122").Value = 5874.7646
$ws.Range("K122").Value = 17624.2938
$ws.Range("M122").Value = -15174.2938
$ws.Range("H132").Value = 2140252.5
$ws.Range("I132").Value = 2749981
$ws.Range("K132").Value = 8249943
$ws.Range("M132").Value = -8247413
$ws.Range("H136").Value = 1863127.5
$ws.Range("I136").Value = 2501730.2
$ws.Range("K136").Value = 7505190.600000001
$ws.Range("M136").Value = -7502640.600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I20").Value = 1798.3334
$ws.Range("J20").Value = 3333
$ws.Range("K20").Value = 1798.3334
$ws.Range("L20").Value = 3333
$ws.Range("M20").Value = -1551.3334
$ws.Range("N20").Value = -3827
$ws.Range("H26").Value = 45000
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("H94").Value = 1853.6086
$ws.Range("I94").Value = 1846.9546
$ws.Range("K94").Value = 1846.9546
$ws.Range("M94").Value = -1395.9546
$ws.Range("H105").Value = 3189
$ws.Range("I105").Value = 3055.0984
$ws.Range("K105").Value = 3055.0984
$ws.Range("M105").Value = -1308.0984
$ws.Range("H134").Value = 3365.1482
$ws.Range("I134").Value = 3042.9524
$ws.Range("K134").Value = 9128.8572
$ws.Range("M134").Value = -6593.8572
$ws.Range("N26").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6144.3335
$ws.Range("I99").Value = 5971.2856
$ws.Range("K99").Value = 5971.2856
$ws.Range("M99").Value = -4473.2856
$ws.Range("H105").Value = 63178.75
$ws.Range("I105").Value = 72061.5
$ws.Range("K105").Value = 72061.5
$ws.Range("M105").Value = -70314.5
$ws.Range("H107").Value = 819.9231
$ws.Range("I107").Value = 568.0625
$ws.Range("J107").Value = 1222.9
$ws.Range("K107").Value = 568.0625
$ws.Range("L107").Value = 1222.9
$ws.Range("M107").Value = 1351.9375
$ws.Range("N107").Value = -5062.9
$ws.Range("H122").Value = 6184
$ws.Range("I122").Value = 5975.6
$ws.Range("K122").Value = 17926.8
$ws.Range("M122").Value = -15476.8
$ws.Range("H126").Value = 6144.3335
$ws.Range("I126").Value = 5971.2856
$ws.Range("K126").Value = 17913.8568
$ws.Range("M126").Value = -15443.8568
$ws.Range("H134").Value = 4056.1562
$ws.Range("I134").Value = 1781.7142
$ws.Range("J134").Value = 5825.1665
$ws.Range("K134").Value = 5345.142599999999
$ws.Range("L134").Value = 17475.4995
$ws.Range("M134").Value = -2810.142599999999
$ws.Range("N134").Value = -22545.4995
$ws.Range("H141").Value = 340173.62
$ws.Range("J141").Value = 361021.4
$ws.Range("L141").Value = 361021.4
$ws.Range("N141").Value = -371381.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 11255132
$ws.Range("J4").Value = 2088413.8
$ws.Range("L4").Value = 6265241.4
$ws.Range("N4").Value = -6265465.4
$ws.Range("H50").Value = 1029.15
$ws.Range("I50").Value = 1034.8
$ws.Range("J50").Value = 1023.5
$ws.Range("K50").Value = 3104.4
$ws.Range("L50").Value = 3070.5
$ws.Range("M50").Value = -2623.4
$ws.Range("N50").Value = -4032.5
$ws.Range("H53").Value = 1029.15
$ws.Range("I53").Value = 1034.8
$ws.Range("J53").Value = 1023.5
$ws.Range("K53").Value = 3104.4
$ws.Range("L53").Value = 3070.5
$ws.Range("M53").Value = -2623.4
$ws.Range("N53").Value = -4032.5
$ws.Range("H114").Value = 260.7143
$ws.Range("I114").Value = 283.83334
$ws.Range("J114").Value = 122
$ws.Range("K114").Value = 851.5000200000001
$ws.Range("L114").Value = 366
$ws.Range("M114").Value = 2402.49998
$ws.Range("N114").Value = -6874
$ws.Range("H117").Value = 4394.5
$ws.Range("I117").Value = 4389
$ws.Range("K117").Value = 13167
$ws.Range("M117").Value = -9725
$ws.Range("H140").Value = 3441.5
$ws.Range("I140").Value = 2467.3333
$ws.Range("J140").Value = 5195
$ws.Range("K140").Value = 7401.999899999999
$ws.Range("L140").Value = 15585
$ws.Range("M140").Value = -2221.999899999999
$ws.Range("N140").Value = -25945

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 724.3333
$ws.Range("I2").Value = 724.3333
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 724.3333
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -611.3333
$ws.Range("H19").Value = 22057.428
$ws.Range("I19").Value = 1752.5
$ws.Range("K19").Value = 1752.5
$ws.Range("M19").Value = -1464.5
$ws.Range("H97").Value = 1407.4706
$ws.Range("I97").Value = 946.4583
$ws.Range("J97").Value = 2513.9
$ws.Range("K97").Value = 946.4583
$ws.Range("L97").Value = 2513.9
$ws.Range("M97").Value = -450.4583
$ws.Range("N97").Value = -3505.9
$ws.Range("H132").Value = 5204.6665
$ws.Range("I132").Value = 5361.3657
$ws.Range("K132").Value = 16084.0971
$ws.Range("M132").Value = -13554.0971
$ws.Range("N2").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 11844.24
$ws.Range("J122").Value = 8513.933999999999
$ws.Range("L122").Value = 25541.802
$ws.Range("N122").Value = -30441.802
$ws.Range("H132").Value = 4184.75
$ws.Range("I132").Value = 4205.08
$ws.Range("J132").Value = 4138.5454
$ws.Range("K132").Value = 12615.24
$ws.Range("L132").Value = 12415.6362
$ws.Range("M132").Value = -10085.24
$ws.Range("N132").Value = -17475.6362
$ws.Range("H136").Value = 6153.577
$ws.Range("I136").Value = 5237.7617
$ws.Range("K136").Value = 15713.2851
$ws.Range("M136").Value = -13163.2851

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 142872180
$ws.Range("I2").Value = 17039.8
$ws.Range("J2").Value = 500010000
$ws.Range("K2").Value = 17039.8
$ws.Range("L2").Value = 500010000
$ws.Range("M2").Value = -16927.8
$ws.Range("N2").Value = -500010224
$ws.Range("H96").Value = 2266.0527
$ws.Range("I96").Value = 1965.8462
$ws.Range("K96").Value = 1965.8462
$ws.Range("M96").Value = -592.8462
$ws.Range("H113").Value = 3623952.8
$ws.Range("J113").Value = 972.6429000000001
$ws.Range("L113").Value = 2917.9287
$ws.Range("N113").Value = -7257.9287
$ws.Range("H122").Value = 3013.7
$ws.Range("I122").Value = 3107.111
$ws.Range("K122").Value = 9321.332999999999
$ws.Range("M122").Value = -6871.332999999999
$ws.Range("H126").Value = 2570.8262
$ws.Range("I126").Value = 2484.7222
$ws.Range("K126").Value = 7454.1666
$ws.Range("M126").Value = -4984.1666
$ws.Range("H136").Value = 16778.4
$ws.Range("I136").Value = 27335.75
$ws.Range("K136").Value = 82007.25
$ws.Range("M136").Value = -79457.25
